$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "SamplesTab" row's query (cell B3) is updated so the Tumor column
# reads directly from samp.sample_tumor_status instead of the collected
# `tumor` list variable.
$newQuery = @'
MATCH (s:study)<--(p:participant)<--(samp:sample)
WHERE s.study_name in ["Discovery of Colorectal Cancer Susceptibility Genes in High-Risk Families"]
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@

$ws.Range("B3").Value = $newQuery

# Move/leave the selection on the cell that was edited.
$ws.Range("B3").Select()
